$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B31: store as a real number (3) instead of text
$ws.Range("B31").Value = 3

# Add new row 32 with annotation data
$ws.Range("A32").Value = "Sunsi Wu"

# B32 keeps the "score as text" quirk from the source data (stored as text "2")
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "2"
$ws.Range("B32").Style = "Normal"

$ws.Range("C32").Value = "not"
$ws.Range("D32").Value = "DFT"
$ws.Range("E32").Value = "MET"
$ws.Range("F32").Value = "aa721c36-81b2-451c-915e-fe15286fe992"
$ws.Range("G32").Value = "SygwwGbRW_annotated.xlsx"
$ws.Range("H32").Value = "This is NOT a proper navigation agent."
